# Insert a new weekly price-report row for "Ajo" (Femacal de La Calera) just
# above the current row 137. This shifts the existing rows 137-210 down to
# 138-211 (growing the used range from A1:R210 to A1:R211) and fills the
# newly opened row 137 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 137..210 down one position.
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new record.
$ws.Cells.Item(137, 1).Value  = 3
$ws.Cells.Item(137, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(137, 3).Value  = "Coquimbo"
$ws.Cells.Item(137, 4).Value  = 44452
$ws.Cells.Item(137, 5).Value  = 5
$ws.Cells.Item(137, 6).Value  = 100112003
$ws.Cells.Item(137, 7).Value  = "Ajo"
$ws.Cells.Item(137, 8).Value  = "Chino"
$ws.Cells.Item(137, 9).Value  = "Primera"
$ws.Cells.Item(137, 10).Value = 83
$ws.Cells.Item(137, 11).Value = 15500
$ws.Cells.Item(137, 12).Value = 16000
$ws.Cells.Item(137, 13).Value = 15771
$ws.Cells.Item(137, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(137, 15).Value = "China"
$ws.Cells.Item(137, 16).Value = 1577
$ws.Cells.Item(137, 17).Value = 10
$ws.Cells.Item(137, 18).Value = "Hortaliza"
